$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9, shifting existing rows 9..78 down to 10..79.
$ws.Rows("9:9").Insert()

# Populate the new row 9 with this week's record. Columns A,B,C,E,F,G,H,I,N,O,Q,R
# are constant across all rows in this sheet; only D,J,K,L,M,P vary per row.
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "Macroferia Regional de Talca"
$ws.Range("C9").Value = "Maule"
$ws.Range("D9").Value = 44552
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 100112030
$ws.Range("G9").Value = "Poroto granado"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 35000
$ws.Range("L9").Value = 35000
$ws.Range("M9").Value = 35000
$ws.Range("N9").Value = "`$/saco 25 kilos"
$ws.Range("O9").Value = "Región del Maule"
$ws.Range("P9").Value = 1400
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"
